$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H header (new shared string "Chỉnh sửa lab 1 nữa")
$ws.Range("H1").Value = "Chỉnh sửa lab 1 nữa"

# New values in column H for several rows
$ws.Range("H2").Value = 65.989999999999995
$ws.Range("H4").Value = 69.010000000000005
$ws.Range("H5").Value = 69.260000000000005
$ws.Range("H11").Value = 69.75
$ws.Range("H12").Value = 69.83

# Updated values in column G
$ws.Range("G11").Value = 68.099999999999994
$ws.Range("G12").Value = 67.900000000000006

# Move the active selection to reflect the new edit location
$ws.Range("H8").Select()
